# Add All Testcases Procedure
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2) text updates ---
$ws.Range("A2").Value = "Testcase steps"
$ws.Range("B2").Value = "Description"

# --- New step descriptions (Read UPCA / REad UPCA_AddOn2) ---
$ws.Range("B12").Value = "Read UPCA"
$ws.Range("B13").Value = "REad UPCA_AddOn2"

# --- Bold header row A2:E2 ---
$ws.Range("A2:E2").Font.Bold = $true

# --- Wrap text for the C/D "result" columns across the data rows ---
$ws.Range("C3:D11").WrapText = $true

# --- Row heights for the rows that now need extra vertical room ---
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 90
$ws.Rows.Item(7).RowHeight = 90
$ws.Rows.Item(10).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 90

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws.Columns.Item(3).ColumnWidth = 62.666666666666664
$ws.Columns.Item(4).ColumnWidth = 60.666666666666664

# --- Update selection to match authored state ---
$ws.Range("E7").Select()
